$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.61535666666667
$ws.Range("H2").Value = 31.84607
$ws.Range("I2").Value = 0.1437966543677284
$ws.Range("J2").Value = 0.1437966543677284
$ws.Range("M2").Value = 14.65767833333333
$ws.Range("N2").Value = 43.973035
$ws.Range("O2").Value = 0.2345581433878666
$ws.Range("P2").Value = 0.2345581433878665
$ws.Range("Q2").Value = 155.5964834136056
$ws.Range("R2").Value = 1400.36835072245
$ws.Range("S2").Value = 0.03372867627388113
$ws.Range("T2").Value = 0.03372867627388113

$ws.Range("G3").Value = 10.61535666666667
$ws.Range("H3").Value = 31.84607
$ws.Range("I3").Value = 0.1437966543677284
$ws.Range("J3").Value = 0.1437966543677284
$ws.Range("M3").Value = 31.695371
$ws.Range("N3").Value = 95.086113
$ws.Range("O3").Value = 0.5072022462686253
$ws.Range("P3").Value = 0.5072022462686253
$ws.Range("Q3").Value = 336.4576678473233
$ws.Range("R3").Value = 3028.11901062591
$ws.Range("S3").Value = 0.07293398610122499
$ws.Range("T3").Value = 0.07293398610122499

$ws.Range("G4").Value = 10.61535666666667
$ws.Range("H4").Value = 31.84607
$ws.Range("I4").Value = 0.1437966543677284
$ws.Range("J4").Value = 0.1437966543677284
$ws.Range("M4").Value = 16.13754733333333
$ws.Range("N4").Value = 48.41264200000001
$ws.Range("O4").Value = 0.2582396103435082
$ws.Range("P4").Value = 0.2582396103435082
$ws.Range("Q4").Value = 171.3058206685489
$ws.Range("R4").Value = 1541.75238601694
$ws.Range("S4").Value = 0.03713399199262232
$ws.Range("T4").Value = 0.03713399199262232

$ws.Range("I5").Value = 0.678584082077682
$ws.Range("J5").Value = 0.678584082077682
$ws.Range("M5").Value = 14.65767833333333
$ws.Range("N5").Value = 43.973035
$ws.Range("O5").Value = 0.2345581433878666
$ws.Range("P5").Value = 0.2345581433878665
$ws.Range("Q5").Value = 734.2681047482895
$ws.Range("R5").Value = 6608.412942734605
$ws.Range("S5").Value = 0.1591674224247007
$ws.Range("T5").Value = 0.1591674224247007

$ws.Range("I6").Value = 0.678584082077682
$ws.Range("J6").Value = 0.678584082077682
$ws.Range("M6").Value = 31.695371
$ws.Range("N6").Value = 95.086113
$ws.Range("O6").Value = 0.5072022462686253
$ws.Range("P6").Value = 0.5072022462686253
$ws.Range("S6").Value = 0.3441793707119335
$ws.Range("T6").Value = 0.3441793707119335

$ws.Range("I7").Value = 0.678584082077682
$ws.Range("J7").Value = 0.678584082077682
$ws.Range("M7").Value = 16.13754733333333
$ws.Range("N7").Value = 48.41264200000001
$ws.Range("O7").Value = 0.2582396103435082
$ws.Range("P7").Value = 0.2582396103435082
$ws.Range("Q7").Value = 808.4013051907251
$ws.Range("R7").Value = 7275.611746716526
$ws.Range("S7").Value = 0.1752372889410478
$ws.Range("T7").Value = 0.1752372889410478

$ws.Range("G8").Value = 13.11220933333333
$ws.Range("H8").Value = 39.336628
$ws.Range("I8").Value = 0.1776192635545896
$ws.Range("J8").Value = 0.1776192635545896
$ws.Range("M8").Value = 14.65767833333333
$ws.Range("N8").Value = 43.973035
$ws.Range("O8").Value = 0.2345581433878666
$ws.Range("P8").Value = 0.2345581433878665
$ws.Range("Q8").Value = 192.1945466473312
$ws.Range("R8").Value = 1729.75091982598
$ws.Range("S8").Value = 0.04166204468928469
$ws.Range("T8").Value = 0.04166204468928469

$ws.Range("G9").Value = 13.11220933333333
$ws.Range("H9").Value = 39.336628
$ws.Range("I9").Value = 0.1776192635545896
$ws.Range("J9").Value = 0.1776192635545896
$ws.Range("M9").Value = 31.695371
$ws.Range("N9").Value = 95.086113
$ws.Range("O9").Value = 0.5072022462686253
$ws.Range("P9").Value = 0.5072022462686253
$ws.Range("Q9").Value = 415.5963394496627
$ws.Range("R9").Value = 3740.367055046964
$ws.Range("S9").Value = 0.09008888945546682
$ws.Range("T9").Value = 0.09008888945546682

$ws.Range("G10").Value = 13.11220933333333
$ws.Range("H10").Value = 39.336628
$ws.Range("I10").Value = 0.1776192635545896
$ws.Range("J10").Value = 0.1776192635545896
$ws.Range("M10").Value = 16.13754733333333
$ws.Range("N10").Value = 48.41264200000001
$ws.Range("O10").Value = 0.2582396103435082
$ws.Range("P10").Value = 0.2582396103435082
$ws.Range("Q10").Value = 211.5988987612418
$ws.Range("R10").Value = 1904.390088851176
$ws.Range("S10").Value = 0.04586832940983811
$ws.Range("T10").Value = 0.04586832940983811

